$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header E1 from "remark" to "Expected "
$ws.Range("E1").Value = "Expected "

# Apply bold font + yellow fill to header row A1:E1 (build format on a helper
# cell first, then copy/paste the format so only a single combined style is
# created instead of two separate intermediate styles)
$helper = $ws.Range("Z100")
$helper.Font.Bold = $true
$helper.Interior.Color = 65535

$headerRange = $ws.Range("A1:E1")
$helper.Copy()
$headerRange.PasteSpecial(-4122)
$helper.Clear()

# Column widths (A, C, D)
$ws.Columns.Item(1).ColumnWidth = 16.667
$ws.Columns.Item(3).ColumnWidth = 11.5
$ws.Columns.Item(4).ColumnWidth = 13.83

# Selection
$ws.Range("B11").Select()

# Page setup orientation (xlPortrait = 1)
$ws.PageSetup.Orientation = 1
